# DistrictBattery type added, EnergyAssets config by defaultEnergyAssets
# - Add DISTRICTBATTERY gridConnection (b7) to config_gridConnections (Table1)
# - Add DistrictBattery defaultEnergyAsset (a3) to config_energyAssets,
#   and repoint the existing PV assets (a1, a2) at the new
#   "Solarpanels_1MW" defaultEnergyAsset reference instead of the literal
#   PHOTOVOLTAIC type2/capacity values.
# - Update sheet selections / active tab to match end-of-edit UI state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. config_gridConnections: add a new DISTRICTBATTERY gridConnection row
# ---------------------------------------------------------------------
$wsGridConn = $wb.Worksheets.Item("config_gridConnections")
$tbl = $wsGridConn.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$wsGridConn.Range("A8").Value = 6
$wsGridConn.Range("B8").Value = "gridConnection"
$wsGridConn.Range("C8").Value = "b7"
$wsGridConn.Range("D8").Value = "DISTRICTBATTERY"
$wsGridConn.Range("E8").Value = "ELECTRIC"
$wsGridConn.Range("F8").Value = "NONE"
$wsGridConn.Range("G8").Value = "NONE"
$wsGridConn.Range("H8").Value = "E2"
$wsGridConn.Range("J8").Value = 1000
$wsGridConn.Range("K8").Value = "hol1"

# ---------------------------------------------------------------------
# 2. config_energyAssets: repoint PV rows at defaultEnergyAssets reference
#    and add the new DistrictBattery defaultEnergyAsset row
# ---------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("config_energyAssets")

$wsAssets.Range("A4").Value = 2
$wsAssets.Range("B4").Value = "energyAsset"
$wsAssets.Range("C4").Value = "a3"
$wsAssets.Range("D4").Value = "STORAGE"

$wsAssets.Range("E2").Value = "Solarpanels_1MW"
$wsAssets.Range("G2").ClearContents()
$wsAssets.Range("H2").ClearContents()

$wsAssets.Range("E3").Value = "Solarpanels_1MW"
$wsAssets.Range("G3").ClearContents()
$wsAssets.Range("H3").ClearContents()

$wsAssets.Range("E4").Value = "District_Battery_500_kWh"
$wsAssets.Range("F4").Value = "b7"

# ---------------------------------------------------------------------
# 3. Update sheet selections to match the final workbook UI state
# ---------------------------------------------------------------------
$wsGridConn.Range("E28").Select()

$wsActors = $wb.Worksheets.Item("config_actors")
$wsActors.Range("D35").Select()

$wsAssets.Activate()
$wsAssets.Range("G1").Select()
